# Apply updated dSF (column F) values to reflect repulled/recalculated data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = -3
$ws.Range("F5").Value = 6
$ws.Range("F8").Value = 6
$ws.Range("F11").Value = 14
$ws.Range("F12").Value = 4
$ws.Range("F15").Value = -3
